$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row($rowA, $rowB) {
    for ($col = 1; $col -le 5; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# Row 2 ("After deadly Taliban attack...") and Row 4 ("Candlelight vigil...") swap places
Swap-Row 2 4

# Row 6 ("Drone strike...") and Row 7 ("Pakistan lifts moratorium...") swap places
Swap-Row 6 7

Write-Output "done"
